$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Salary" column header
$ws.Range("D1").Value = "Salary"

# Fix second row's LastName from "Jannette" to "Jackson"
$ws.Range("B2").Value = "Jackson"

# Populate Salary values for each employee
$ws.Range("D2").Value = 200000.0
$ws.Range("D3").Value = 110000.0
$ws.Range("D4").Value = 135000.0
$ws.Range("D5").Value = 125000.0
